# Update the F-column (time_taken) timestamps on the "data" sheet to reflect
# a later re-query of the PanelApp API, and add a new "metadata" sheet that
# records details about that query (data_name/id/version/... panel_query_time).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$newTimestamps = @(
  "2021-10-05 14:33:12.562206",
  "2021-10-05 14:33:12.562213",
  "2021-10-05 14:33:12.562216",
  "2021-10-05 14:33:12.562218",
  "2021-10-05 14:33:12.562220",
  "2021-10-05 14:33:12.562222",
  "2021-10-05 14:33:12.562225",
  "2021-10-05 14:33:12.562227",
  "2021-10-05 14:33:12.562229",
  "2021-10-05 14:33:12.562232",
  "2021-10-05 14:33:12.562234",
  "2021-10-05 14:33:12.562236",
  "2021-10-05 14:33:12.562238",
  "2021-10-05 14:33:12.562240",
  "2021-10-05 14:33:12.562242",
  "2021-10-05 14:33:12.562245",
  "2021-10-05 14:33:12.562247",
  "2021-10-05 14:33:12.562249",
  "2021-10-05 14:33:12.562251",
  "2021-10-05 14:33:12.562254",
  "2021-10-05 14:33:12.562256",
  "2021-10-05 14:33:12.562258",
  "2021-10-05 14:33:12.562260",
  "2021-10-05 14:33:12.562262",
  "2021-10-05 14:33:12.562265",
  "2021-10-05 14:33:12.562267",
  "2021-10-05 14:33:12.562269",
  "2021-10-05 14:33:12.562272",
  "2021-10-05 14:33:12.562274",
  "2021-10-05 14:33:12.562276",
  "2021-10-05 14:33:12.562278",
  "2021-10-05 14:33:12.562280",
  "2021-10-05 14:33:12.562282",
  "2021-10-05 14:33:12.562285",
  "2021-10-05 14:33:12.562287",
  "2021-10-05 14:33:12.562289",
  "2021-10-05 14:33:12.562291",
  "2021-10-05 14:33:12.562293",
  "2021-10-05 14:33:12.562295",
  "2021-10-05 14:33:12.562297",
  "2021-10-05 14:33:12.562300",
  "2021-10-05 14:33:12.562302",
  "2021-10-05 14:33:12.562304",
  "2021-10-05 14:33:12.562307",
  "2021-10-05 14:33:12.562309",
  "2021-10-05 14:33:12.562311",
  "2021-10-05 14:33:12.562313",
  "2021-10-05 14:33:12.562316",
  "2021-10-05 14:33:12.562318",
  "2021-10-05 14:33:12.562320",
  "2021-10-05 14:33:12.562322",
  "2021-10-05 14:33:12.562324",
  "2021-10-05 14:33:12.562326",
  "2021-10-05 14:33:12.562329",
  "2021-10-05 14:33:12.562331",
  "2021-10-05 14:33:12.562333",
  "2021-10-05 14:33:12.562335",
  "2021-10-05 14:33:12.562337",
  "2021-10-05 14:33:12.562339",
  "2021-10-05 14:33:12.562341",
  "2021-10-05 14:33:12.562343",
  "2021-10-05 14:33:12.562346",
  "2021-10-05 14:33:12.562348",
  "2021-10-05 14:33:12.562350",
  "2021-10-05 14:33:12.562353",
  "2021-10-05 14:33:12.562356",
  "2021-10-05 14:33:12.562358",
  "2021-10-05 14:33:12.562360",
  "2021-10-05 14:33:12.562363",
  "2021-10-05 14:33:12.562365",
  "2021-10-05 14:33:12.562367",
  "2021-10-05 14:33:12.562369",
  "2021-10-05 14:33:12.562371",
  "2021-10-05 14:33:12.562374",
  "2021-10-05 14:33:12.562376",
  "2021-10-05 14:33:12.562378",
  "2021-10-05 14:33:12.562382",
  "2021-10-05 14:33:12.562384",
  "2021-10-05 14:33:12.562387",
  "2021-10-05 14:33:12.562389",
  "2021-10-05 14:33:12.562391",
  "2021-10-05 14:33:12.562393",
  "2021-10-05 14:33:12.562395",
  "2021-10-05 14:33:12.562397",
  "2021-10-05 14:33:12.562400",
  "2021-10-05 14:33:12.562402",
  "2021-10-05 14:33:12.562404",
  "2021-10-05 14:33:12.562406",
  "2021-10-05 14:33:12.562408",
  "2021-10-05 14:33:12.562410",
  "2021-10-05 14:33:12.562412",
  "2021-10-05 14:33:12.562414",
  "2021-10-05 14:33:12.562417",
  "2021-10-05 14:33:12.562420",
  "2021-10-05 14:33:12.562422",
  "2021-10-05 14:33:12.562424",
  "2021-10-05 14:33:12.562426",
  "2021-10-05 14:33:12.562428",
  "2021-10-05 14:33:12.562430",
  "2021-10-05 14:33:12.562432",
  "2021-10-05 14:33:12.562435",
  "2021-10-05 14:33:12.562437",
  "2021-10-05 14:33:12.562439",
  "2021-10-05 14:33:12.562441",
  "2021-10-05 14:33:12.562443",
  "2021-10-05 14:33:12.562445",
  "2021-10-05 14:33:12.562447",
  "2021-10-05 14:33:12.562449",
  "2021-10-05 14:33:12.562453",
  "2021-10-05 14:33:12.562456",
  "2021-10-05 14:33:12.562458",
  "2021-10-05 14:33:12.562460",
  "2021-10-05 14:33:12.562462",
  "2021-10-05 14:33:12.562464",
  "2021-10-05 14:33:12.562467",
  "2021-10-05 14:33:12.562469",
  "2021-10-05 14:33:12.562471",
  "2021-10-05 14:33:12.562473",
  "2021-10-05 14:33:12.562475",
  "2021-10-05 14:33:12.562477",
  "2021-10-05 14:33:12.562479",
  "2021-10-05 14:33:12.562481",
  "2021-10-05 14:33:12.562484",
  "2021-10-05 14:33:12.562486",
  "2021-10-05 14:33:12.562488",
  "2021-10-05 14:33:12.562490",
  "2021-10-05 14:33:12.562492",
  "2021-10-05 14:33:12.562494",
  "2021-10-05 14:33:12.562498",
  "2021-10-05 14:33:12.562501",
  "2021-10-05 14:33:12.562503",
  "2021-10-05 14:33:12.562505",
  "2021-10-05 14:33:12.562507",
  "2021-10-05 14:33:12.562509",
  "2021-10-05 14:33:12.562512",
  "2021-10-05 14:33:12.562514",
  "2021-10-05 14:33:12.562516",
  "2021-10-05 14:33:12.562518",
  "2021-10-05 14:33:12.562520",
  "2021-10-05 14:33:12.562522",
  "2021-10-05 14:33:12.562524",
  "2021-10-05 14:33:12.562526",
  "2021-10-05 14:33:12.562528",
  "2021-10-05 14:33:12.562530",
  "2021-10-05 14:33:12.562533",
  "2021-10-05 14:33:12.562535",
  "2021-10-05 14:33:12.562537",
  "2021-10-05 14:33:12.562539",
  "2021-10-05 14:33:12.562542",
  "2021-10-05 14:33:12.562544",
  "2021-10-05 14:33:12.562546",
  "2021-10-05 14:33:12.562548",
  "2021-10-05 14:33:12.562550",
  "2021-10-05 14:33:12.562552",
  "2021-10-05 14:33:12.562554",
  "2021-10-05 14:33:12.562557",
  "2021-10-05 14:33:12.562559",
  "2021-10-05 14:33:12.562561",
  "2021-10-05 14:33:12.562563",
  "2021-10-05 14:33:12.562565",
  "2021-10-05 14:33:12.562567",
  "2021-10-05 14:33:12.562569",
  "2021-10-05 14:33:12.562571",
  "2021-10-05 14:33:12.562573",
  "2021-10-05 14:33:12.562576",
  "2021-10-05 14:33:12.562578",
  "2021-10-05 14:33:12.562580",
  "2021-10-05 14:33:12.562582",
  "2021-10-05 14:33:12.562584",
  "2021-10-05 14:33:12.562586",
  "2021-10-05 14:33:12.562588",
  "2021-10-05 14:33:12.562590",
  "2021-10-05 14:33:12.562594",
  "2021-10-05 14:33:12.562597",
  "2021-10-05 14:33:12.562599",
  "2021-10-05 14:33:12.562601",
  "2021-10-05 14:33:12.562603",
  "2021-10-05 14:33:12.562605",
  "2021-10-05 14:33:12.562607",
  "2021-10-05 14:33:12.562609",
  "2021-10-05 14:33:12.562611",
  "2021-10-05 14:33:12.562613",
  "2021-10-05 14:33:12.562615",
  "2021-10-05 14:33:12.562617",
  "2021-10-05 14:33:12.562619",
  "2021-10-05 14:33:12.562622",
  "2021-10-05 14:33:12.562624",
  "2021-10-05 14:33:12.562626",
  "2021-10-05 14:33:12.562628",
  "2021-10-05 14:33:12.562630",
  "2021-10-05 14:33:12.562632",
  "2021-10-05 14:33:12.562634",
  "2021-10-05 14:33:12.562636",
  "2021-10-05 14:33:12.562639",
  "2021-10-05 14:33:12.562641",
  "2021-10-05 14:33:12.562643",
  "2021-10-05 14:33:12.562645",
  "2021-10-05 14:33:12.562647",
  "2021-10-05 14:33:12.562650",
  "2021-10-05 14:33:12.562652",
  "2021-10-05 14:33:12.562654",
  "2021-10-05 14:33:12.562656",
  "2021-10-05 14:33:12.562658",
  "2021-10-05 14:33:12.562660",
  "2021-10-05 14:33:12.562663",
  "2021-10-05 14:33:12.562665",
  "2021-10-05 14:33:12.562667",
  "2021-10-05 14:33:12.562669",
  "2021-10-05 14:33:12.562671",
  "2021-10-05 14:33:12.562673",
  "2021-10-05 14:33:12.562675",
  "2021-10-05 14:33:12.562677",
  "2021-10-05 14:33:12.562680",
  "2021-10-05 14:33:12.562682",
  "2021-10-05 14:33:12.562684",
  "2021-10-05 14:33:12.562686",
  "2021-10-05 14:33:12.562689",
  "2021-10-05 14:33:12.562691",
  "2021-10-05 14:33:12.562693",
  "2021-10-05 14:33:12.562695",
  "2021-10-05 14:33:12.562697",
  "2021-10-05 14:33:12.562699",
  "2021-10-05 14:33:12.562701",
  "2021-10-05 14:33:12.562703",
  "2021-10-05 14:33:12.562706"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# Add the new "metadata" worksheet right after "data"
$meta = $wb.Worksheets.Add($null, $ws1)
$meta.Name = "metadata"

# Header row (B1:G1) - reuse the bold/bordered header style from the data sheet
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
$cols = @("B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $meta.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
    $ws1.Range("B1").Copy()
    $cell.PasteSpecial(-4122)
}

# Row 2 - the metadata values for this panel
$a2 = $meta.Range("A2")
$a2.Value = 0
$ws1.Range("A2").Copy()
$a2.PasteSpecial(-4122)

$meta.Range("B2").Value = "Ataxia - paediatric"
$meta.Range("C2").Value = 271

$d2 = $meta.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "0.293"
$d2.ClearFormats()

$meta.Range("E2").Value = "2021-10-04T07:11:15.290045Z"
$meta.Range("F2").Value = "2021-10-05 14:33:12.559317"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/271/?format=json"

# Keep "data" as the active/selected sheet (adding a sheet makes it active by default)
$ws1.Activate()

Write-Output "metadata sheet added and time_taken column refreshed"
